# The post previously stored in row 605 ("「幸せは伝染性。喜びを広めよう」")
# was removed from the source data. Delete that entire row so every
# subsequent row shifts up by one (old row 606 becomes 605, ... old row
# 712 becomes 711), matching the new extent A1:C711.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(605).Delete()
